# Append two new data rows (134, 135) to Sheet1, mirroring the existing
# "date / site / 24 hourly values" pattern already present in the sheet
# (e.g. rows 132/133), then move the active selection the way Excel does
# after the user finishes entering data in the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$siteFangSquarePing = "四方坪站充电量(kw)"
$siteGaoLing = "高岭站充电量(kw)"

# New row 134: 2025-11-06 (serial 45967), 四方坪站
$ws.Cells.Item(134, 1).Value = 45967
$ws.Cells.Item(134, 2).Value = $siteFangSquarePing

$row134 = @(546.71900000000005, 960.6389999999999, 382.25800000000004, 416.26900000000001, 382.01600000000002, 497.27700000000004, 338.74299999999999, 164.655, 116.19800000000002, 203.88099999999997, 176.322, 210.744, 949.0780000000002, 1139.5529999999999, 372.87, 522.47799999999995, 236.14099999999999, 152.84100000000001, 23.16, 94.579999999999984, 41.28, 53.28, 58.929999999999993, 63.64)

$col = 3
foreach ($v in $row134) {
    $ws.Cells.Item(134, $col).Value = $v
    $col = $col + 1
}

# New row 135: 2025-11-06 (serial 45967), 高岭站
$ws.Cells.Item(135, 1).Value = 45967
$ws.Cells.Item(135, 2).Value = $siteGaoLing

$row135 = @(403.916, 192.77699999999999, 39.664000000000001, 95.661000000000001, 54.812000000000005, 232.149, 324.27700000000004, 100.821, 90.173000000000002, 116.274, 228.197, 158.10499999999999, 390.56899999999996, 431.11300000000006, 418.12399999999997, 27.13, 142.16399999999999, 21.695999999999998, 109.58699999999999, 144.59299999999999, 131.41300000000001, 66.688000000000002, 0, 20.922000000000001)

$col = 3
foreach ($v in $row135) {
    $ws.Cells.Item(135, $col).Value = $v
    $col = $col + 1
}

# Update the active cell / selection like Excel would after typing into
# the new rows (activeCell S137, same cell selected).
$ws.Range("S137").Select() | Out-Null
